$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2400
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 2400
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("N65").Value = -26240
$ws.Range("H111").Value = 10874.866
$ws.Range("I111").Value = 25737.4
$ws.Range("J111").Value = 3443.6
$ws.Range("K111").Value = 77212.20000000001
$ws.Range("L111").Value = 10330.8
$ws.Range("M111").Value = -74145.20000000001
$ws.Range("N111").Value = -16464.8
$ws.Range("H129").Value = 400658.84
$ws.Range("J129").Value = 421880.6
$ws.Range("L129").Value = 1265641.8
$ws.Range("N129").Value = -1275641.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 78678.16
$ws.Range("I2").Value = 1792.3
$ws.Range("J2").Value = 334964.34
$ws.Range("K2").Value = 1792.3
$ws.Range("L2").Value = 334964.34
$ws.Range("M2").Value = -1679.3
$ws.Range("N2").Value = -335190.34
$ws.Range("H32").Value = 28718.254
$ws.Range("I32").Value = 4485.852
$ws.Range("J32").Value = 174112.67
$ws.Range("K32").Value = 4485.852
$ws.Range("L32").Value = 174112.67
$ws.Range("M32").Value = -4198.852
$ws.Range("N32").Value = -174686.67
$ws.Range("H45").Value = 38859.703
$ws.Range("I45").Value = 72831.64
$ws.Range("K45").Value = 72831.64
$ws.Range("M45").Value = -72454.64
$ws.Range("H106").Value = 45980
$ws.Range("J106").Value = 45980
$ws.Range("L106").Value = 45980
$ws.Range("N106").Value = -48504
$ws.Range("H116").Value = 78678.16
$ws.Range("I116").Value = 1792.3
$ws.Range("J116").Value = 334964.34
$ws.Range("K116").Value = 1792.3
$ws.Range("L116").Value = 334964.34
$ws.Range("M116").Value = 501.7
$ws.Range("N116").Value = -339552.34
$ws.Range("H122").Value = 3031.5557
$ws.Range("I122").Value = 2550.6667
$ws.Range("J122").Value = 3993.3333
$ws.Range("K122").Value = 7652.000100000001
$ws.Range("L122").Value = 11979.9999
$ws.Range("M122").Value = -5202.000100000001
$ws.Range("N122").Value = -16879.9999
$ws.Range("H132").Value = 2934.5386
$ws.Range("I132").Value = 3019.375
$ws.Range("J132").Value = 2798.8
$ws.Range("K132").Value = 9058.125
$ws.Range("L132").Value = 8396.400000000001
$ws.Range("M132").Value = -6528.125
$ws.Range("N132").Value = -13456.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 78678.16
$ws.Range("I3").Value = 1792.3
$ws.Range("J3").Value = 334964.34
$ws.Range("K3").Value = 1792.3
$ws.Range("L3").Value = 334964.34
$ws.Range("M3").Value = -1678.3
$ws.Range("N3").Value = -335192.34
$ws.Range("H64").Value = 593.4
$ws.Range("I64").Value = 90
$ws.Range("J64").Value = 649.3333
$ws.Range("K64").Value = 90
$ws.Range("L64").Value = 649.3333
$ws.Range("M64").Value = 135
$ws.Range("N64").Value = -1099.3333
$ws.Range("H67").Value = 593.4
$ws.Range("I67").Value = 90
$ws.Range("J67").Value = 649.3333
$ws.Range("K67").Value = 90
$ws.Range("L67").Value = 649.3333
$ws.Range("M67").Value = 690
$ws.Range("N67").Value = -2209.3333
$ws.Range("H132").Value = 68000
$ws.Range("J132").Value = 68000
$ws.Range("L132").Value = 68000
$ws.Range("N132").Value = -78120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4318
$ws.Range("I122").Value = 382.72726
$ws.Range("J122").Value = 25962
$ws.Range("K122").Value = 3444.54534
$ws.Range("L122").Value = 233658
$ws.Range("M122").Value = -994.5453400000001
$ws.Range("N122").Value = -238558
$ws.Range("H123").Value = 4118.25
$ws.Range("I123").Value = 2486.6667
$ws.Range("J123").Value = 4662.1113
$ws.Range("K123").Value = 7460.000100000001
$ws.Range("L123").Value = 13986.3339
$ws.Range("M123").Value = -5010.000100000001
$ws.Range("N123").Value = -18886.3339
$ws.Range("H124").Value = 3155.111
$ws.Range("I124").Value = 2999.5
$ws.Range("J124").Value = 3174.5625
$ws.Range("K124").Value = 8998.5
$ws.Range("L124").Value = 9523.6875
$ws.Range("M124").Value = -4088.5
$ws.Range("N124").Value = -19343.6875
$ws.Range("H125").Value = 800
$ws.Range("I125").Value = 800
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 2400
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = 2520
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 1257.5
$ws.Range("I126").Value = 1010
$ws.Range("K126").Value = 3030
$ws.Range("M126").Value = 1910
$ws.Range("H129").Value = 11598264
$ws.Range("I129").Value = 35722732
$ws.Range("J129").Value = 340177.6
$ws.Range("K129").Value = 107168196
$ws.Range("L129").Value = 1020532.8
$ws.Range("M129").Value = -107163196
$ws.Range("N129").Value = -1030532.8
$ws.Range("H131").Value = 715.04
$ws.Range("J131").Value = 775.03656
$ws.Range("L131").Value = 2325.10968
$ws.Range("N131").Value = -12405.10968

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 275121.97
$ws.Range("I102").Value = 2386.9333
$ws.Range("K102").Value = 2386.9333
$ws.Range("M102").Value = -764.9333000000001
$ws.Range("H113").Value = 3003.1428
$ws.Range("I113").Value = 3170.3333
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 3170.3333
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -1000.3333
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 1133.3334
$ws.Range("I122").Value = 950
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 2850
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -400
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 3239.9546
$ws.Range("I132").Value = 2686.8667
$ws.Range("J132").Value = 4425.143
$ws.Range("K132").Value = 8060.6001
$ws.Range("L132").Value = 13275.429
$ws.Range("M132").Value = -5530.6001
$ws.Range("N132").Value = -18335.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 49371.81
$ws.Range("J40").Value = 2031.0769
$ws.Range("L40").Value = 2031.0769
$ws.Range("N40").Value = -2303.0769
$ws.Range("H44").Value = 30673
$ws.Range("J44").Value = 30673
$ws.Range("L44").Value = 30673
$ws.Range("N44").Value = -31585
$ws.Range("H115").Value = 20302
$ws.Range("J115").Value = 20302
$ws.Range("L115").Value = 20302
$ws.Range("N115").Value = -22652
$ws.Range("H122").Value = 1498.8
$ws.Range("I122").Value = 1498.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4496.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2046.4
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 21315
$ws.Range("I2").Value = 14000
$ws.Range("J2").Value = 32287.5
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 32287.5
$ws.Range("M2").Value = -13888
$ws.Range("N2").Value = -32511.5
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 2185.611
$ws.Range("I122").Value = 1148.909
$ws.Range("J122").Value = 3814.7144
$ws.Range("K122").Value = 3446.727
$ws.Range("L122").Value = 11444.1432
$ws.Range("M122").Value = -996.7270000000003
$ws.Range("N122").Value = -16344.1432
$ws.Range("H132").Value = 1287.6666
$ws.Range("I132").Value = 835.7222
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 2507.1666
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = 22.83339999999998
$ws.Range("N132").Value = -17057.9999
$ws.Range("H136").Value = 1150
$ws.Range("I136").Value = 866.6667
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2600.0001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -50.0001000000002
$ws.Range("N136").Value = -11100

Write-Output "applied changes"